$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd URL in B25 (missing leading "h" in "https")
$ws.Range("B25").Value = "https://www.google.com/maps/place/Heidekreis-Klinikum+GmbH+Krankenhaus+Soltau/@52.9894409,9.847291,15z/data=!3m1!4b1!4m8!1m2!11m1!2s1tsS4C8icZfBtXgqho9ekuv3aB34!3m4!1s0x47b1b"

# Add hyperlink for the fixed cell
$ws.Hyperlinks.Add($ws.Range("B25"), "https://www.google.com/maps/place/Heidekreis-Klinikum+GmbH+Krankenhaus+Soltau/@52.9894409,9.847291,15z/data=!3m1!4b1!4m8!1m2!11m1!2s1tsS4C8icZfBtXgqho9ekuv3aB34!3m4!1s0x47b1b")

# Widen column B (target stored width 37.6640625 characters)
$ws.Columns("B").ColumnWidth = 36.83

# Move selection to B26
$ws.Range("B26").Select()
